# Sort the curvature-data rows (everything below the header row) by
# column A (time) in ascending order, leaving the header row untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstDataRow = $used.Row + 1                 # header is the first row -> data starts on the next row
$lastDataRow  = $used.Row + $used.Rows.Count - 1
$firstCol     = $used.Column
$numCols      = $used.Columns.Count

# Read all data rows into an array of row-arrays
$rows = @()
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @()
    for ($i = 0; $i -lt $numCols; $i++) {
        $c = $firstCol + $i
        $rowVals += , $ws.Cells.Item($r, $c).Value2
    }
    $rows += , $rowVals
}

# Sort the rows by the first column (time) ascending
$sortedRows = $rows | Sort-Object { $_[0] }

# Write the sorted rows back in place
$r = $firstDataRow
foreach ($rowVals in $sortedRows) {
    for ($i = 0; $i -lt $numCols; $i++) {
        $c = $firstCol + $i
        $ws.Cells.Item($r, $c).Value2 = $rowVals[$i]
    }
    $r++
}
